# Auto-generated edit script: updates LeveProfit-sheet price/profit values
# per the authoritative diff (scheduled price-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 807968.1  # H17: 678769.1 -> 807968.1
$ws.Cells.Item(17, 10).Value = 807968.1  # J17: 678769.1 -> 807968.1
$ws.Cells.Item(17, 12).Value = 2423904.3  # L17: 2036307.3 -> 2423904.3
$ws.Cells.Item(17, 14).Value = -2424240.3  # N17: -2036643.3 -> -2424240.3
$ws.Cells.Item(33, 8).Value = 437.77274  # H33: 513.7 -> 437.77274
$ws.Cells.Item(33, 9).Value = 224.1579  # I33: 189.56 -> 224.1579
$ws.Cells.Item(33, 10).Value = 1790.6666  # J33: 2134.4 -> 1790.6666
$ws.Cells.Item(33, 11).Value = 224.1579  # K33: 189.56 -> 224.1579
$ws.Cells.Item(33, 12).Value = 1790.6666  # L33: 2134.4 -> 1790.6666
$ws.Cells.Item(33, 13).Value = 4.842099999999988  # M33: 39.44 -> 4.842099999999988
$ws.Cells.Item(33, 14).Value = -2248.6666  # N33: -2592.4 -> -2248.6666
$ws.Cells.Item(109, 8).Value = 25690  # H109: 29000 -> 25690
$ws.Cells.Item(109, 10).Value = 25690  # J109: 29000 -> 25690
$ws.Cells.Item(109, 12).Value = 25690  # L109: 29000 -> 25690
$ws.Cells.Item(109, 14).Value = -28464  # N109: -31774 -> -28464
$ws.Cells.Item(132, 8).Value = 138491.14  # H132: 111425.18 -> 138491.14
$ws.Cells.Item(132, 9).Value = 1209.803  # I132: 1092.7808 -> 1209.803
$ws.Cells.Item(132, 10).Value = 1432858  # J132: 558884.3 -> 1432858
$ws.Cells.Item(132, 11).Value = 3629.409000000001  # K132: 3278.3424 -> 3629.409000000001
$ws.Cells.Item(132, 12).Value = 4298574  # L132: 1676652.9 -> 4298574
$ws.Cells.Item(132, 13).Value = -1099.409000000001  # M132: -748.3424 -> -1099.409000000001
$ws.Cells.Item(132, 14).Value = -4303634  # N132: -1681712.9 -> -4303634
$ws.Cells.Item(137, 8).Value = 59994.707  # H137: 42718.332 -> 59994.707
$ws.Cells.Item(137, 9).Value = 125923.75  # I137: 67514.664 -> 125923.75
$ws.Cells.Item(137, 11).Value = 377771.25  # K137: 202543.992 -> 377771.25
$ws.Cells.Item(137, 13).Value = -375221.25  # M137: -199993.992 -> -375221.25
$ws.Cells.Item(138, 8).Value = 6038072.5  # H138: 6508500 -> 6038072.5
$ws.Cells.Item(138, 9).Value = 1583.5  # I138: 1838.4117 -> 1583.5
$ws.Cells.Item(138, 10).Value = 8215167  # J138: 8352054 -> 8215167
$ws.Cells.Item(138, 11).Value = 4750.5  # K138: 5515.2351 -> 4750.5
$ws.Cells.Item(138, 12).Value = 24645501  # L138: 25056162 -> 24645501
$ws.Cells.Item(138, 13).Value = 389.5  # M138: -375.2350999999999 -> 389.5
$ws.Cells.Item(138, 14).Value = -24655781  # N138: -25066442 -> -24655781

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 162951.5  # H32: 165346.38 -> 162951.5
$ws.Cells.Item(32, 9).Value = 167314.95  # I32: 170101.86 -> 167314.95
$ws.Cells.Item(32, 11).Value = 167314.95  # K32: 170101.86 -> 167314.95
$ws.Cells.Item(32, 13).Value = -167027.95  # M32: -169814.86 -> -167027.95
$ws.Cells.Item(61, 8).Value = 892.8570999999999  # H61: 905.4583 -> 892.8570999999999
$ws.Cells.Item(61, 9).Value = 810  # I61: 830.0952 -> 810
$ws.Cells.Item(61, 10).Value = 1100  # J61: 1433 -> 1100
$ws.Cells.Item(61, 11).Value = 810  # K61: 830.0952 -> 810
$ws.Cells.Item(61, 12).Value = 1100  # L61: 1433 -> 1100
$ws.Cells.Item(61, 13).Value = -598  # M61: -618.0952 -> -598
$ws.Cells.Item(61, 14).Value = -1524  # N61: -1857 -> -1524
$ws.Cells.Item(74, 8).Value = 68868.97  # H74: 50447.953 -> 68868.97
$ws.Cells.Item(74, 9).Value = 81968.36  # I74: 60344.03 -> 81968.36
$ws.Cells.Item(74, 10).Value = 3372  # J74: 2381.2856 -> 3372
$ws.Cells.Item(74, 11).Value = 81968.36  # K74: 60344.03 -> 81968.36
$ws.Cells.Item(74, 12).Value = 3372  # L74: 2381.2856 -> 3372
$ws.Cells.Item(74, 13).Value = -81094.36  # M74: -59470.03 -> -81094.36
$ws.Cells.Item(74, 14).Value = -5120  # N74: -4129.2856 -> -5120
$ws.Cells.Item(77, 8).Value = 68868.97  # H77: 50447.953 -> 68868.97
$ws.Cells.Item(77, 9).Value = 81968.36  # I77: 60344.03 -> 81968.36
$ws.Cells.Item(77, 10).Value = 3372  # J77: 2381.2856 -> 3372
$ws.Cells.Item(77, 11).Value = 409841.8  # K77: 301720.15 -> 409841.8
$ws.Cells.Item(77, 12).Value = 16860  # L77: 11906.428 -> 16860
$ws.Cells.Item(77, 13).Value = -405473.8  # M77: -297352.15 -> -405473.8
$ws.Cells.Item(77, 14).Value = -25596  # N77: -20642.428 -> -25596
$ws.Cells.Item(122, 8).Value = 900.93335  # H122: 812.3 -> 900.93335
$ws.Cells.Item(122, 9).Value = 728.5454999999999  # I122: 770.8421 -> 728.5454999999999
$ws.Cells.Item(122, 10).Value = 1375  # J122: 1600 -> 1375
$ws.Cells.Item(122, 11).Value = 2185.6365  # K122: 2312.5263 -> 2185.6365
$ws.Cells.Item(122, 12).Value = 4125  # L122: 4800 -> 4125
$ws.Cells.Item(122, 13).Value = 264.3635000000004  # M122: 137.4737 -> 264.3635000000004
$ws.Cells.Item(122, 14).Value = -9025  # N122: -9700 -> -9025
$ws.Cells.Item(136, 8).Value = 892.8570999999999  # H136: 905.4583 -> 892.8570999999999
$ws.Cells.Item(136, 9).Value = 810  # I136: 830.0952 -> 810
$ws.Cells.Item(136, 10).Value = 1100  # J136: 1433 -> 1100
$ws.Cells.Item(136, 11).Value = 2430  # K136: 2490.2856 -> 2430
$ws.Cells.Item(136, 12).Value = 3300  # L136: 4299 -> 3300
$ws.Cells.Item(136, 13).Value = 120  # M136: 59.71439999999984 -> 120
$ws.Cells.Item(136, 14).Value = -8400  # N136: -9399 -> -8400

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 50915.477  # H134: 28387.475 -> 50915.477
$ws.Cells.Item(134, 9).Value = 1880.7333  # I134: 1115.9375 -> 1880.7333
$ws.Cells.Item(134, 10).Value = 173502.33  # J134: 173835.67 -> 173502.33
$ws.Cells.Item(134, 11).Value = 5642.199900000001  # K134: 3347.8125 -> 5642.199900000001
$ws.Cells.Item(134, 12).Value = 520506.99  # L134: 521507.01 -> 520506.99
$ws.Cells.Item(134, 13).Value = -3107.199900000001  # M134: -812.8125 -> -3107.199900000001
$ws.Cells.Item(134, 14).Value = -525576.99  # N134: -526577.01 -> -525576.99

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 26289.24  # H31: 28901.918 -> 26289.24
$ws.Cells.Item(31, 9).Value = 31770.77  # I31: 35323.715 -> 31770.77
$ws.Cells.Item(31, 10).Value = 12037.267  # J31: 12847.429 -> 12037.267
$ws.Cells.Item(31, 11).Value = 31770.77  # K31: 35323.715 -> 31770.77
$ws.Cells.Item(31, 12).Value = 12037.267  # L31: 12847.429 -> 12037.267
$ws.Cells.Item(31, 13).Value = -31475.77  # M31: -35028.715 -> -31475.77
$ws.Cells.Item(31, 14).Value = -12627.267  # N31: -13437.429 -> -12627.267
$ws.Cells.Item(34, 8).Value = 26289.24  # H34: 28901.918 -> 26289.24
$ws.Cells.Item(34, 9).Value = 31770.77  # I34: 35323.715 -> 31770.77
$ws.Cells.Item(34, 10).Value = 12037.267  # J34: 12847.429 -> 12037.267
$ws.Cells.Item(34, 11).Value = 31770.77  # K34: 35323.715 -> 31770.77
$ws.Cells.Item(34, 12).Value = 12037.267  # L34: 12847.429 -> 12037.267
$ws.Cells.Item(34, 13).Value = -31568.77  # M34: -35121.715 -> -31568.77
$ws.Cells.Item(34, 14).Value = -12441.267  # N34: -13251.429 -> -12441.267
$ws.Cells.Item(62, 8).Value = 3037.2222  # H62: 3004.1667 -> 3037.2222
$ws.Cells.Item(62, 9).Value = 3005  # I62: 2959.0908 -> 3005
$ws.Cells.Item(62, 10).Value = 3150  # J62: 3500 -> 3150
$ws.Cells.Item(62, 11).Value = 3005  # K62: 2959.0908 -> 3005
$ws.Cells.Item(62, 12).Value = 3150  # L62: 3500 -> 3150
$ws.Cells.Item(62, 13).Value = -2381  # M62: -2335.0908 -> -2381
$ws.Cells.Item(62, 14).Value = -4398  # N62: -4748 -> -4398
$ws.Cells.Item(65, 8).Value = 3037.2222  # H65: 3004.1667 -> 3037.2222
$ws.Cells.Item(65, 9).Value = 3005  # I65: 2959.0908 -> 3005
$ws.Cells.Item(65, 10).Value = 3150  # J65: 3500 -> 3150
$ws.Cells.Item(65, 11).Value = 15025  # K65: 14795.454 -> 15025
$ws.Cells.Item(65, 12).Value = 15750  # L65: 17500 -> 15750
$ws.Cells.Item(65, 13).Value = -11905  # M65: -11675.454 -> -11905
$ws.Cells.Item(65, 14).Value = -21990  # N65: -23740 -> -21990
$ws.Cells.Item(106, 8).Value = 58333.332  # H106: 56223.668 -> 58333.332
$ws.Cells.Item(106, 10).Value = 58333.332  # J106: 56223.668 -> 58333.332
$ws.Cells.Item(106, 12).Value = 58333.332  # L106: 56223.668 -> 58333.332
$ws.Cells.Item(106, 14).Value = -60857.332  # N106: -58747.668 -> -60857.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 6421.6665  # H5: 8545.691999999999 -> 6421.6665
$ws.Cells.Item(5, 9).Value = 899.1667  # I5: 899.1429000000001 -> 899.1667
$ws.Cells.Item(5, 11).Value = 2697.5001  # K5: 2697.4287 -> 2697.5001
$ws.Cells.Item(5, 13).Value = -2585.5001  # M5: -2585.4287 -> -2585.5001
$ws.Cells.Item(26, 8).Value = 400000500  # H26: 250000300 -> 400000500
$ws.Cells.Item(26, 9).Value = 1000000000  # I26: 500000060 -> 1000000000
$ws.Cells.Item(26, 10).Value = 833  # J26: 549.5 -> 833
$ws.Cells.Item(26, 11).Value = 3000000000  # K26: 1500000180 -> 3000000000
$ws.Cells.Item(26, 12).Value = 2499  # L26: 1648.5 -> 2499
$ws.Cells.Item(26, 13).Value = -2999999712  # M26: -1499999892 -> -2999999712
$ws.Cells.Item(26, 14).Value = -3075  # N26: -2224.5 -> -3075
$ws.Cells.Item(42, 8).Value = 55559390  # H42: 83335590 -> 55559390
$ws.Cells.Item(42, 10).Value = 55559390  # J42: 83335590 -> 55559390
$ws.Cells.Item(42, 12).Value = 166678170  # L42: 250006770 -> 166678170
$ws.Cells.Item(42, 14).Value = -166679238  # N42: -250007838 -> -166679238
$ws.Cells.Item(103, 8).Value = 1531.25  # H103: 2113.8462 -> 1531.25
$ws.Cells.Item(103, 9).Value = 708.3333  # I103: 1062.5 -> 708.3333
$ws.Cells.Item(103, 10).Value = 4000  # J103: 2581.111 -> 4000
$ws.Cells.Item(103, 11).Value = 2124.9999  # K103: 3187.5 -> 2124.9999
$ws.Cells.Item(103, 12).Value = 12000  # L103: 7743.333 -> 12000
$ws.Cells.Item(103, 13).Value = -1245.9999  # M103: -2308.5 -> -1245.9999
$ws.Cells.Item(103, 14).Value = -13758  # N103: -9501.332999999999 -> -13758
$ws.Cells.Item(106, 8).Value = 4499.846  # H106: 4472.5293 -> 4499.846
$ws.Cells.Item(106, 10).Value = 4499.846  # J106: 4472.5293 -> 4499.846
$ws.Cells.Item(106, 12).Value = 13499.538  # L106: 13417.5879 -> 13499.538
$ws.Cells.Item(106, 14).Value = -15391.538  # N106: -15309.5879 -> -15391.538
$ws.Cells.Item(113, 8).Value = 663.2105  # H113: 672.8333 -> 663.2105
$ws.Cells.Item(113, 10).Value = 675.6667  # J113: 686.58826 -> 675.6667
$ws.Cells.Item(113, 12).Value = 2027.0001  # L113: 2059.76478 -> 2027.0001
$ws.Cells.Item(113, 14).Value = -6367.0001  # N113: -6399.76478 -> -6367.0001
$ws.Cells.Item(131, 8).Value = 23585706  # H131: 20492620 -> 23585706
$ws.Cells.Item(131, 9).Value = 498.18182  # I131: 510 -> 498.18182
$ws.Cells.Item(131, 10).Value = 29762784  # J131: 24510680 -> 29762784
$ws.Cells.Item(131, 11).Value = 1494.54546  # K131: 1530 -> 1494.54546
$ws.Cells.Item(131, 12).Value = 89288352  # L131: 73532040 -> 89288352
$ws.Cells.Item(131, 13).Value = 3545.45454  # M131: 3510 -> 3545.45454
$ws.Cells.Item(131, 14).Value = -89298432  # N131: -73542120 -> -89298432
$ws.Cells.Item(134, 8).Value = 4117.3716  # H134: 4111.643 -> 4117.3716
$ws.Cells.Item(134, 9).Value = 1200.8572  # I134: 1199.5 -> 1200.8572
$ws.Cells.Item(134, 10).Value = 8492.143  # J134: 7315 -> 8492.143
$ws.Cells.Item(134, 11).Value = 3602.5716  # K134: 3598.5 -> 3602.5716
$ws.Cells.Item(134, 12).Value = 25476.429  # L134: 21945 -> 25476.429
$ws.Cells.Item(134, 13).Value = 1467.4284  # M134: 1471.5 -> 1467.4284
$ws.Cells.Item(134, 14).Value = -35616.429  # N134: -32085 -> -35616.429
$ws.Cells.Item(135, 8).Value = 6421.6665  # H135: 8545.691999999999 -> 6421.6665
$ws.Cells.Item(135, 9).Value = 899.1667  # I135: 899.1429000000001 -> 899.1667
$ws.Cells.Item(135, 11).Value = 8092.5003  # K135: 8092.2861 -> 8092.5003
$ws.Cells.Item(135, 13).Value = -5557.5003  # M135: -5557.2861 -> -5557.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 20314.13  # H102: 20008.61 -> 20314.13
$ws.Cells.Item(102, 9).Value = 11560.615  # I102: 9062.706 -> 11560.615
$ws.Cells.Item(102, 10).Value = 31693.7  # J102: 51022 -> 31693.7
$ws.Cells.Item(102, 11).Value = 11560.615  # K102: 9062.706 -> 11560.615
$ws.Cells.Item(102, 12).Value = 31693.7  # L102: 51022 -> 31693.7
$ws.Cells.Item(102, 13).Value = -9938.615  # M102: -7440.706 -> -9938.615
$ws.Cells.Item(102, 14).Value = -34937.7  # N102: -54266 -> -34937.7
$ws.Cells.Item(126, 8).Value = 1520.0952  # H126: 1600.8 -> 1520.0952
$ws.Cells.Item(126, 9).Value = 1386.375  # I126: 1445 -> 1386.375
$ws.Cells.Item(126, 10).Value = 1948  # J126: 1799.091 -> 1948
$ws.Cells.Item(126, 11).Value = 4159.125  # K126: 4335 -> 4159.125
$ws.Cells.Item(126, 12).Value = 5844  # L126: 5397.272999999999 -> 5844
$ws.Cells.Item(126, 13).Value = -1689.125  # M126: -1865 -> -1689.125
$ws.Cells.Item(126, 14).Value = -10784  # N126: -10337.273 -> -10784
$ws.Cells.Item(127, 8).Value = 55326  # H127: 0 -> 55326
$ws.Cells.Item(127, 10).Value = 55326  # J127: 0 -> 55326
$ws.Cells.Item(127, 12).Value = 55326  # L127: 0 -> 55326
$ws.Cells.Item(127, 14).Value = -65246  # N127: None -> -65246
$ws.Cells.Item(130, 8).Value = 0  # H130: 25756 -> 0
$ws.Cells.Item(130, 10).Value = 0  # J130: 25756 -> 0
$ws.Cells.Item(130, 12).Value = 0  # L130: 25756 -> 0
$ws.Cells.Item(130, 14).Value = ""  # N130: clear (was -35796)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(104, 8).Value = 29185  # H104: 17013.334 -> 29185
$ws.Cells.Item(104, 10).Value = 29185  # J104: 17013.334 -> 29185
$ws.Cells.Item(104, 12).Value = 29185  # L104: 17013.334 -> 29185
$ws.Cells.Item(104, 14).Value = -36173  # N104: -24001.334 -> -36173
$ws.Cells.Item(138, 8).Value = 40000  # H138: 42000 -> 40000
$ws.Cells.Item(138, 10).Value = 40000  # J138: 42000 -> 40000
$ws.Cells.Item(138, 12).Value = 40000  # L138: 42000 -> 40000
$ws.Cells.Item(138, 14).Value = -50280  # N138: -52280 -> -50280

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(104, 8).Value = 27408  # H104: 0 -> 27408
$ws.Cells.Item(104, 10).Value = 27408  # J104: 0 -> 27408
$ws.Cells.Item(104, 12).Value = 27408  # L104: 0 -> 27408
$ws.Cells.Item(104, 14).Value = -34396  # N104: None -> -34396
$ws.Cells.Item(136, 8).Value = 359080.2  # H136: 913138.2 -> 359080.2
$ws.Cells.Item(136, 9).Value = 1793.0435  # I136: 4573.5713 -> 1793.0435
$ws.Cells.Item(136, 10).Value = 2002601  # J136: 2503126.2 -> 2002601
$ws.Cells.Item(136, 11).Value = 5379.1305  # K136: 13720.7139 -> 5379.1305
$ws.Cells.Item(136, 12).Value = 6007803  # L136: 7509378.600000001 -> 6007803
$ws.Cells.Item(136, 13).Value = -2829.1305  # M136: -11170.7139 -> -2829.1305
$ws.Cells.Item(136, 14).Value = -6012903  # N136: -7514478.600000001 -> -6012903
